# Apply updates to column D ("value") for the specified rows on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    @{Row=10;  New=1},
    @{Row=31;  New=0},
    @{Row=40;  New=1},
    @{Row=50;  New=0},
    @{Row=61;  New=1},
    @{Row=77;  New=0},
    @{Row=121; New=1},
    @{Row=131; New=0},
    @{Row=140; New=4},
    @{Row=161; New=2},
    @{Row=177; New=0},
    @{Row=179; New=21},
    @{Row=210; New=13},
    @{Row=221; New=0},
    @{Row=231; New=1},
    @{Row=251; New=3},
    @{Row=261; New=38},
    @{Row=267; New=2},
    @{Row=277; New=4},
    @{Row=310; New=2},
    @{Row=360; New=1},
    @{Row=377; New=0},
    @{Row=410; New=8},
    @{Row=451; New=12},
    @{Row=461; New=3},
    @{Row=467; New=1},
    @{Row=477; New=2},
    @{Row=500; New=5},
    @{Row=510; New=6},
    @{Row=550; New=0},
    @{Row=561; New=1},
    @{Row=577; New=3},
    @{Row=581; New=0},
    @{Row=600; New=1}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 4).Value = $u.New
}
